# AH Up in 20200921
# - WEB_UI!D2  : fill in the "Real-RESULT" so it matches the expected result in E2
# - WIN_UI!E2  : fill in the "Real-RESULT" so it matches the expected result in F2

$wb = $excel.ActiveWorkbook

$web = $wb.Worksheets.Item("WEB_UI")
$web.Range("D2").Value = "Allen_百度搜索"

$win = $wb.Worksheets.Item("WIN_UI")
$win.Range("E2").Value = "显示为 15"
